$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = 44336
$ws.Cells.Item(2, 11).Value = 'Winter Nelis'
$ws.Cells.Item(2, 13).Value = 250
$ws.Cells.Item(2, 14).Value = 21000
$ws.Cells.Item(2, 15).Value = 22000
$ws.Cells.Item(2, 16).Value = 21500
$ws.Cells.Item(2, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(2, 19).Value = 1194
$ws.Cells.Item(3, 4).Value = 44398
$ws.Cells.Item(3, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(3, 13).Value = 200
$ws.Cells.Item(3, 14).Value = 20000
$ws.Cells.Item(3, 15).Value = 21000
$ws.Cells.Item(3, 16).Value = 20500
$ws.Cells.Item(3, 17).Value = '$/caja 20 kilos empedrada'
$ws.Cells.Item(3, 19).Value = 1025
$ws.Cells.Item(3, 20).Value = 20
$ws.Cells.Item(4, 4).Value = 44398
$ws.Cells.Item(4, 11).Value = 'Winter Nelis'
$ws.Cells.Item(4, 13).Value = 200
$ws.Cells.Item(4, 14).Value = 20000
$ws.Cells.Item(4, 15).Value = 21000
$ws.Cells.Item(4, 16).Value = 20500
$ws.Cells.Item(4, 17).Value = '$/caja 20 kilos empedrada'
$ws.Cells.Item(4, 19).Value = 1025
$ws.Cells.Item(4, 20).Value = 20
$ws.Cells.Item(5, 4).Value = 44642
$ws.Cells.Item(5, 12).Value = 'Segunda'
$ws.Cells.Item(5, 13).Value = 270
$ws.Cells.Item(5, 14).Value = 19000
$ws.Cells.Item(5, 15).Value = 20000
$ws.Cells.Item(5, 16).Value = 19500
$ws.Cells.Item(5, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(5, 19).Value = 1083
$ws.Cells.Item(6, 4).Value = 44474
$ws.Cells.Item(6, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(6, 12).Value = 'Segunda'
$ws.Cells.Item(6, 13).Value = 270
$ws.Cells.Item(6, 14).Value = 18000
$ws.Cells.Item(6, 15).Value = 19000
$ws.Cells.Item(6, 16).Value = 18500
$ws.Cells.Item(6, 17).Value = '$/caja 18 kilos empedrada'
$ws.Cells.Item(6, 19).Value = 1028
$ws.Cells.Item(7, 4).Value = 44474
$ws.Cells.Item(7, 11).Value = 'Winter Nelis'
$ws.Cells.Item(7, 13).Value = 250
$ws.Cells.Item(8, 4).Value = 44313
$ws.Cells.Item(8, 12).Value = 'Tercera'
$ws.Cells.Item(8, 14).Value = 15000
$ws.Cells.Item(8, 15).Value = 16000
$ws.Cells.Item(8, 16).Value = 15500
$ws.Cells.Item(8, 19).Value = 861
$ws.Cells.Item(9, 4).Value = 44525
$ws.Cells.Item(9, 13).Value = 300
$ws.Cells.Item(9, 14).Value = 19000
$ws.Cells.Item(9, 16).Value = 19500
$ws.Cells.Item(9, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(9, 19).Value = 1083
$ws.Cells.Item(9, 20).Value = 18
$ws.Cells.Item(10, 4).Value = 44497
$ws.Cells.Item(10, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(10, 13).Value = 300
$ws.Cells.Item(10, 14).Value = 17000
$ws.Cells.Item(10, 15).Value = 18000
$ws.Cells.Item(10, 16).Value = 17500
$ws.Cells.Item(10, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(10, 19).Value = 972
$ws.Cells.Item(10, 20).Value = 18
$ws.Cells.Item(11, 4).Value = 44497
$ws.Cells.Item(11, 11).Value = 'Winter Nelis'
$ws.Cells.Item(11, 13).Value = 250
$ws.Cells.Item(11, 14).Value = 17000
$ws.Cells.Item(11, 15).Value = 18000
$ws.Cells.Item(11, 16).Value = 17500
$ws.Cells.Item(11, 19).Value = 972
$ws.Cells.Item(12, 4).Value = 44314
$ws.Cells.Item(12, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(12, 14).Value = 17000
$ws.Cells.Item(12, 15).Value = 18000
$ws.Cells.Item(12, 16).Value = 17500
$ws.Cells.Item(12, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(12, 19).Value = 972
$ws.Cells.Item(13, 4).Value = 44329
$ws.Cells.Item(13, 13).Value = 340
$ws.Cells.Item(13, 14).Value = 21000
$ws.Cells.Item(13, 16).Value = 21500
$ws.Cells.Item(13, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(13, 19).Value = 1194
$ws.Cells.Item(14, 4).Value = 44323
$ws.Cells.Item(14, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(14, 14).Value = 15000
$ws.Cells.Item(14, 15).Value = 16000
$ws.Cells.Item(14, 16).Value = 15500
$ws.Cells.Item(14, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(14, 19).Value = 861
$ws.Cells.Item(15, 4).Value = 44421
$ws.Cells.Item(15, 13).Value = 270
$ws.Cells.Item(15, 14).Value = 16000
$ws.Cells.Item(15, 16).Value = 16500
$ws.Cells.Item(15, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(15, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(15, 19).Value = 917
$ws.Cells.Item(16, 4).Value = 44421
$ws.Cells.Item(16, 14).Value = 16000
$ws.Cells.Item(16, 16).Value = 16500
$ws.Cells.Item(16, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(16, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(16, 19).Value = 917
$ws.Cells.Item(17, 4).Value = 44292
$ws.Cells.Item(17, 14).Value = 22000
$ws.Cells.Item(17, 15).Value = 23000
$ws.Cells.Item(17, 16).Value = 22500
$ws.Cells.Item(17, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(17, 19).Value = 1250
$ws.Cells.Item(18, 4).Value = 44292
$ws.Cells.Item(18, 11).Value = 'Winter Nelis'
$ws.Cells.Item(18, 13).Value = 250
$ws.Cells.Item(18, 14).Value = 22000
$ws.Cells.Item(18, 15).Value = 23000
$ws.Cells.Item(18, 16).Value = 22500
$ws.Cells.Item(18, 19).Value = 1250
$ws.Cells.Item(19, 4).Value = 44341
$ws.Cells.Item(19, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(19, 13).Value = 300
$ws.Cells.Item(19, 14).Value = 17000
$ws.Cells.Item(19, 15).Value = 18000
$ws.Cells.Item(19, 16).Value = 17500
$ws.Cells.Item(19, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(19, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(19, 19).Value = 972
$ws.Cells.Item(20, 4).Value = 44280
$ws.Cells.Item(20, 13).Value = 350
$ws.Cells.Item(20, 14).Value = 24000
$ws.Cells.Item(20, 15).Value = 25000
$ws.Cells.Item(20, 16).Value = 24500
$ws.Cells.Item(20, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(20, 19).Value = 1361
$ws.Cells.Item(21, 4).Value = 44280
$ws.Cells.Item(21, 11).Value = 'Winter Nelis'
$ws.Cells.Item(21, 13).Value = 300
$ws.Cells.Item(21, 14).Value = 24000
$ws.Cells.Item(21, 15).Value = 25000
$ws.Cells.Item(21, 16).Value = 24500
$ws.Cells.Item(21, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(21, 19).Value = 1361
$ws.Cells.Item(22, 4).Value = 44355
$ws.Cells.Item(22, 13).Value = 200
$ws.Cells.Item(22, 14).Value = 17000
$ws.Cells.Item(22, 16).Value = 17500
$ws.Cells.Item(22, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(22, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(22, 19).Value = 972
$ws.Cells.Item(23, 4).Value = 44355
$ws.Cells.Item(23, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(23, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(24, 4).Value = 44371
$ws.Cells.Item(24, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(24, 12).Value = 'Calibre 90'
$ws.Cells.Item(24, 13).Value = 140
$ws.Cells.Item(24, 14).Value = 17000
$ws.Cells.Item(24, 15).Value = 18000
$ws.Cells.Item(24, 16).Value = 17429
$ws.Cells.Item(24, 17).Value = '$/caja 18 kilos embalada'
$ws.Cells.Item(24, 19).Value = 968
$ws.Cells.Item(25, 4).Value = 44371
$ws.Cells.Item(25, 11).Value = 'Winter Nelis'
$ws.Cells.Item(25, 12).Value = 'Calibre 80'
$ws.Cells.Item(25, 13).Value = 120
$ws.Cells.Item(25, 14).Value = 17000
$ws.Cells.Item(25, 15).Value = 18000
$ws.Cells.Item(25, 16).Value = 17500
$ws.Cells.Item(25, 17).Value = '$/caja 18 kilos embalada'
$ws.Cells.Item(25, 19).Value = 972
